$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.06592966666667
$ws.Range("H2").Value = 99.197789
$ws.Range("I2").Value = 0.4620579289161133
$ws.Range("J2").Value = 0.4620579289161132
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 162.98837
$ws.Range("N2").Value = 488.96511
$ws.Range("O2").Value = 0.9909539753179891
$ws.Range("P2").Value = 0.9909539753179891
$ws.Range("Q2").Value = 5389.361978904644
$ws.Range("R2").Value = 48504.25781014179
$ws.Range("S2").Value = 0.4578781414866193
$ws.Range("T2").Value = 0.4578781414866192
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.06592966666667
$ws.Range("H3").Value = 99.197789
$ws.Range("I3").Value = 0.4620579289161133
$ws.Range("J3").Value = 0.4620579289161132
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6513563333333333
$ws.Range("N3").Value = 1.954069
$ws.Range("O3").Value = 0.003960185305646138
$ws.Range("P3").Value = 0.003960185305646138
$ws.Range("Q3").Value = 21.53770270593789
$ws.Range("R3").Value = 193.839324353441
$ws.Range("S3").Value = 0.00182983502045088
$ws.Range("T3").Value = 0.001829835020450879
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.06592966666667
$ws.Range("H4").Value = 99.197789
$ws.Range("I4").Value = 0.4620579289161133
$ws.Range("J4").Value = 0.4620579289161132
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8364996666666666
$ws.Range("N4").Value = 2.509499
$ws.Range("O4").Value = 0.005085839376364744
$ws.Range("P4").Value = 0.005085839376364744
$ws.Range("Q4").Value = 27.65963914419011
$ws.Range("R4").Value = 248.936752297711
$ws.Range("S4").Value = 0.002349952409043111
$ws.Range("T4").Value = 0.00234995240904311
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.69090766666667
$ws.Range("H5").Value = 44.072723
$ws.Range("I5").Value = 0.2052883568914378
$ws.Range("J5").Value = 0.2052883568914378
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 162.98837
$ws.Range("N5").Value = 488.96511
$ws.Range("O5").Value = 0.9909539753179891
$ws.Range("P5").Value = 0.9909539753179891
$ws.Range("Q5").Value = 2394.447094410503
$ws.Range("R5").Value = 21550.02384969453
$ws.Range("S5").Value = 0.2034313133480684
$ws.Range("T5").Value = 0.2034313133480684
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.69090766666667
$ws.Range("H6").Value = 44.072723
$ws.Range("I6").Value = 0.2052883568914378
$ws.Range("J6").Value = 0.2052883568914378
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6513563333333333
$ws.Range("N6").Value = 1.954069
$ws.Range("O6").Value = 0.003960185305646138
$ws.Range("P6").Value = 0.003960185305646138
$ws.Range("Q6").Value = 9.569015751098554
$ws.Range("R6").Value = 86.121141759887
$ws.Range("S6").Value = 0.0008129799343817123
$ws.Range("T6").Value = 0.0008129799343817123
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.69090766666667
$ws.Range("H7").Value = 44.072723
$ws.Range("I7").Value = 0.2052883568914378
$ws.Range("J7").Value = 0.2052883568914378
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8364996666666666
$ws.Range("N7").Value = 2.509499
$ws.Range("O7").Value = 0.005085839376364744
$ws.Range("P7").Value = 0.005085839376364744
$ws.Range("Q7").Value = 12.28893936619744
$ws.Range("R7").Value = 110.600454295777
$ws.Range("S7").Value = 0.001044063608987693
$ws.Range("T7").Value = 0.001044063608987693
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.805466
$ws.Range("H8").Value = 71.416398
$ws.Range("I8").Value = 0.3326537141924489
$ws.Range("J8").Value = 0.3326537141924489
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 162.98837
$ws.Range("N8").Value = 488.96511
$ws.Range("O8").Value = 0.9909539753179891
$ws.Range("P8").Value = 0.9909539753179891
$ws.Range("Q8").Value = 3880.01410043042
$ws.Range("R8").Value = 34920.12690387378
$ws.Range("S8").Value = 0.3296445204833014
$ws.Range("T8").Value = 0.3296445204833014
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.805466
$ws.Range("H9").Value = 71.416398
$ws.Range("I9").Value = 0.3326537141924489
$ws.Range("J9").Value = 0.3326537141924489
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6513563333333333
$ws.Range("N9").Value = 1.954069
$ws.Range("O9").Value = 0.003960185305646138
$ws.Range("P9").Value = 0.003960185305646138
$ws.Range("Q9").Value = 15.50584104705133
$ws.Range("R9").Value = 139.552569423462
$ws.Range("S9").Value = 0.001317370350813546
$ws.Range("T9").Value = 0.001317370350813546
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.805466
$ws.Range("H10").Value = 71.416398
$ws.Range("I10").Value = 0.3326537141924489
$ws.Range("J10").Value = 0.3326537141924489
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8364996666666666
$ws.Range("N10").Value = 2.509499
$ws.Range("O10").Value = 0.005085839376364744
$ws.Range("P10").Value = 0.005085839376364744
$ws.Range("Q10").Value = 19.91326437384467
$ws.Range("R10").Value = 179.219379364602
$ws.Range("S10").Value = 0.00169182335833394
$ws.Range("T10").Value = 0.00169182335833394

Write-Host "Applied NATMI Sema6d-Kdr updates"